$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 6.051899999999988
$ws.Range("A9").Value = -20.34129999999998
$ws.Range("A18").Value = -23.03110000000002
$ws.Range("A20").Value = -22.11180000000002
$ws.Range("E21").Value = 13.22869999999999
